$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8617137670516968
$ws.Range("B1").Value = 3.110463857650757
$ws.Range("C1").Value = 3.006486177444458
$ws.Range("D1").Value = 1.707107663154602
$ws.Range("E1").Value = 1.311923027038574
